$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rev. C -> Rev. D : rename the sheet (cascades into the scoped Print_Area
# defined names' sheet-qualified references).
$ws.Name = "Domino LED Rev. D"

# The legacy "bare" _xlnm.Print_Area entry doesn't get its RefersTo string
# rewritten by the rename above, so fix it up explicitly.
$wb.Names.Item(1).RefersTo = "='Domino LED Rev. D'!`$A`$1:`$I`$10"

# Mirror the extra accumulated Print_Area_* history entry that shows up
# after this revision (another sheet-scoped Print_Area copy).
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0", "='Domino LED Rev. D'!`$A`$1:`$I`$1")

# Silkscreen font ratio change (20%) nudged the rendered glyph widths, which
# LibreOffice/Excel re-measured as a small (~0.5%) bump in every column's
# character-width based on the BOM table columns (A:I).
$ws.Range("A1").EntireColumn.ColumnWidth = 4
$ws.Range("B1").EntireColumn.ColumnWidth = 4
$ws.Range("C1").EntireColumn.ColumnWidth = 25.8333333333333
$ws.Range("D1").EntireColumn.ColumnWidth = 27.6666666666667
$ws.Range("E1").EntireColumn.ColumnWidth = 30.8333333333333
$ws.Range("F1").EntireColumn.ColumnWidth = 25.8333333333333
$ws.Range("G1").EntireColumn.ColumnWidth = 40
$ws.Range("H1").EntireColumn.ColumnWidth = 62
$ws.Range("I1").EntireColumn.ColumnWidth = 24
